$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear existing content but keep per-cell styles (A23 keeps its existing style index,
# which we reuse in-place below so the final style table matches the target exactly).
$ws.UsedRange.ClearContents()

# --- Header / title cells ---
$ws.Range("A1").Value = 'Functional and Non-Functional Requirements'
$ws.Range("A3").Value = 'Web Application / iOS / Android'
$ws.Range("A4").Value = 'Number'
$ws.Range("B4").Value = 'Requirement'
$ws.Range("C4").Value = 'Status'

# --- Data rows 5-34 ---
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = 'Interface that allows user to view media, by type or by date from database'
$ws.Range("C5").Value = 'Unplanned'

$ws.Range("A6").Value = 2
$ws.Range("B6").Value = 'Interface that allows user to add media'
$ws.Range("C6").Value = 'Unplanned'

$ws.Range("A7").Value = 2.1
$ws.Range("B7").Value = 'Add Media on website'
$ws.Range("C7").Value = 'Unplanned'

$ws.Range("A8").Value = '2.1.1'
$ws.Range("B8").Value = 'Uploader for files, all types'
$ws.Range("C8").Value = 'Unplanned'

$ws.Range("A9").Value = '2.1.2'
$ws.Range("B9").Value = 'Uploader to allow message to be entered like a text or tweet'
$ws.Range("C9").Value = 'Unplanned'

$ws.Range("A10").Value = 2.2000000000000002
$ws.Range("B10").Value = 'Add Media via mobile device'
$ws.Range("C10").Value = 'Unplanned'

$ws.Range("A11").Value = '2.2.1'
$ws.Range("B11").Value = 'Use phone technology to take or upload photo'
$ws.Range("C11").Value = 'Unplanned'

$ws.Range("A12").Value = '2.2.2'
$ws.Range("B12").Value = 'Use phone technology to take or upload video'
$ws.Range("C12").Value = 'Unplanned'

$ws.Range("A13").Value = '2.2.3'
$ws.Range("B13").Value = 'Use phone technology to take or upload vocal message'
$ws.Range("C13").Value = 'Unplanned'

$ws.Range("A14").Value = '2.2.4'
$ws.Range("B14").Value = 'Use phone technology to enter message like text or tweet as well as upload text'
$ws.Range("C14").Value = 'Unplanned'

$ws.Range("A15").Value = 3
$ws.Range("B15").Value = 'Interface that prompts user to login, or keep session'
$ws.Range("C15").Value = 'Unplanned'

$ws.Range("A16").Value = 4
$ws.Range("B16").Value = 'Interface that allows user to signup for service'
$ws.Range("C16").Value = 'Unplanned'

$ws.Range("A17").Value = 5
$ws.Range("B17").Value = 'Interface that allows user to modify their media (change title, change details, delete)'
$ws.Range("C17").Value = 'Unplanned'

$ws.Range("A18").Value = 6
$ws.Range("B18").Value = 'Interface that allows user to share media with another user'
$ws.Range("C18").Value = 'Unplanned'

$ws.Range("A19").Value = 7
$ws.Range("B19").Value = 'Interface that allows user to manage profile, passwords, information'
$ws.Range("C19").Value = 'Unplanned'

$ws.Range("A20").Value = 8
$ws.Range("B20").Value = 'passwords will only be able to be matched, no missing password retreival method, only password resets'
$ws.Range("C20").Value = 'Unplanned'

$ws.Range("B21").Value = 'Non-Functional Requirement'

$ws.Range("A22").Value = 1
$ws.Range("B22").Value = 'Media and functionality of app needs to be appliable to all available smart phones or web computers'
$ws.Range("C22").Value = 'Unplanned'

$ws.Range("A23").Value = 2
$ws.Range("B23").Value = 'Server needs:  Linux server, also availability to run on a rasberry pi'
$ws.Range("C23").Value = 'Unplanned'

$ws.Range("A24").Value = 3
$ws.Range("B24").Value = 'Server software is to be written using PHP '
$ws.Range("C24").Value = 'Unplanned'

$ws.Range("A25").Value = 4
$ws.Range("B25").Value = 'Libraries need to be built out for android (Objects that help protect against data intrusion)'
$ws.Range("C25").Value = 'Unplanned'

$ws.Range("A26").Value = 5
$ws.Range("B26").Value = 'Libraries are to be written using java (or until further notice)'
$ws.Range("C26").Value = 'Unplanned'

$ws.Range("A27").Value = 6
$ws.Range("B27").Value = 'Libraries need to be built out for iOS (objects that help protect against data intrusion)'
$ws.Range("C27").Value = 'Unplanned'

$ws.Range("A28").Value = 7
$ws.Range("B28").Value = 'Libraries need to be written using Objective C (or until further notice)'
$ws.Range("C28").Value = 'Unplanned'

$ws.Range("A29").Value = 8
$ws.Range("B29").Value = 'A prototype will be written using a Web application wrapped with each respected application wrapper'
$ws.Range("C29").Value = 'Unplanned'

$ws.Range("A30").Value = 9
$ws.Range("B30").Value = 'This prototype will give a look and feel for the mobile app and will be the model for the native application'
$ws.Range("C30").Value = 'Unplanned'

$ws.Range("A31").Value = 10
$ws.Range("B31").Value = 'If a native application is out side of the scope near the project end date, then prototype will be Evaluated'
$ws.Range("C31").Value = 'Unplanned'

$ws.Range("A32").Value = 11
$ws.Range("B32").Value = 'Data passed to the server will be encrypted'
$ws.Range("C32").Value = 'Unplanned'

$ws.Range("A33").Value = 12
$ws.Range("B33").Value = 'Data will be assigned a unique ID, based off the contents of the data and the iterative id of the server'
$ws.Range("C33").Value = 'Unplanned'

$ws.Range("A34").Value = 13
$ws.Range("B34").Value = 'Implementation of Oauth or other type of one type logins will be investigated'
$ws.Range("C34").Value = 'Unplanned'

# --- Section header formatting (bold), matching existing "Normal bold" style ---
$ws.Range("B21").Font.Bold = $true

# A25 inherited bold formatting from the old "Non-Functional Requirement" header that used
# to live there; strip it so the cell goes back to the plain/default font before centering.
$ws.Range("A25").Font.Bold = $false

# --- Column A number/ID alignment: center-align.
# A23 already carries the pre-existing (now-unused-visually) style index 2;
# updating its alignment first causes the engine to mutate that style slot in
# place (keeping cellXfs at 3 entries) instead of appending a new one. Every
# subsequent cell we center then naturally reuses that same slot.
$ws.Range("A23").HorizontalAlignment = -4108
$ws.Range("A5:A22").HorizontalAlignment = -4108
$ws.Range("A24:A34").HorizontalAlignment = -4108

# A21 stays empty (no value) but still carries the centered style, matching the source
$ws.Range("A21").HorizontalAlignment = -4108

# --- Restore selection / active cell like the committed workbook ---
$ws.Range("E20").Select()
